$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 650 (the "彼は完璧ではない..." post), shifting all following rows up by one.
$ws.Rows.Item(650).Delete()
